# Update Mappings 22 Ontologies
# Adds a new "MOP_DEF" column (F) with definition text for each BFO class row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell F1: "MOP_DEF" ---
# Copy the formatting from E1 (the neighbouring header cell) so F1 picks up
# the same bold/centered/bordered header style, then set its text.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "MOP_DEF"

# --- Data cells F2:F4: definition text for each row ---
$ws.Range("F2").Value = "['p is a process if p is an occurrent that has temporal proper parts and for some time t, p specifically depends on some material entity at t. [BFO]', locstr(`"Process, i.e., a physical entity with a temporal evolution that 'has a meaning for the ontologist'`", 'en')]"

$ws.Range("F3").Value = "['B is a process profile if there is some process c such that b is process profile of c. B is process profile of c holds when b is a proper occurrent part of c and there is some proper occurrent part d of c which has no parts in common with b and is mutually dependent on b and is such that b, c and d occupy the same temporal region. [BFO]']"

$ws.Range("F4").Value = "['B is a disposition means: b is a realizable entity and b’s bearer is some material entity and b is such that if it ceases to exist, then its bearer is physically changed, and b’s realization occurs when and because this bearer is in some special physical circumstances, and this realization occurs in virtue of the bearer’s physical make-up. [BFO]']"
